$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column remains stored as text, matching the source data
# format (values like "1.024" or "318.14" must not be coerced into numeric
# cells by Excel's automatic type detection). The Volume(1h) column values
# already contain spaces/percent signs so they stay text without this.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.499.84"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "1.874.64"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "1.026"
$ws.Range("E4").Value = "  +2.80%  "
$ws.Range("D5").Value = "318.14"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "1.022"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "0.5146"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").Value = "0.3971"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("D9").Value = "0.08369"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "1.113"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "42.18"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("D12").Value = "6.257"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "20.54"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").Value = "1.027"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.236"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.809.92"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "0.00001111"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "91.31"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "0.06780"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "1.022"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Value = "5.971"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "28.532.14"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "2.287"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "162.33"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").Value = "2.034.52"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "20.85"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").Value = "2.366"
$ws.Range("D30").Value = "127.79"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "0.1053"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "1.040"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").Value = "5.821"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").Value = "3.649"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").Value = "0.02431"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "0.06499"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "0.2189"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "1.279"
$ws.Range("E38").Value = "  +5.07%  "
$ws.Range("D39").Value = "8.883"
$ws.Range("E39").Value = "  -6.36%  "
$ws.Range("D40").Value = "1.189"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "0.6436"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "5.038"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "11.25"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "0.6038"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "13.04"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "3.733"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("D47").Value = "1.229"
$ws.Range("E47").Value = "  -3.65%  "
$ws.Range("D48").Value = "1.993"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "1.210"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "122.09"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "0.06864"
$ws.Range("E51").Value = "  +0.00%  "
